# Insert a new weekly price record as row 53 ("Fruta / hortaliza, semanal"),
# pushing the existing rows 53-110 down to 54-111.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(53).Insert()

$ws.Range("A53").Value = 10
$ws.Range("B53").Value = "Vega Modelo de Temuco"
$ws.Range("C53").Value = "La Araucanía"
$ws.Range("D53").Value = 44539
$ws.Range("E53").Value = 9
$ws.Range("F53").Value = 100112012
$ws.Range("G53").Value = "Espinaca"
$ws.Range("H53").Value = "Sin especificar"
$ws.Range("I53").Value = "Primera"
$ws.Range("J53").Value = 95
$ws.Range("K53").Value = 7000
$ws.Range("L53").Value = 8000
$ws.Range("M53").Value = 7368
$ws.Range("N53").Value = "$/docena de atados"
$ws.Range("O53").Value = "Región de La Araucanía"
$ws.Range("P53").Value = 2456
$ws.Range("Q53").Value = 3
$ws.Range("R53").Value = "Hortaliza"
